$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H43").Value = 1383.3334
$ws.Range("I43").Value = 1200
$ws.Range("J43").Value = 1750
$ws.Range("K43").Value = 1200
$ws.Range("L43").Value = 1750
$ws.Range("M43").Value = -1131
$ws.Range("N43").Value = -1888

$ws.Range("H106").Value = 10443.5
$ws.Range("I106").Value = 6162.222
$ws.Range("J106").Value = 18149.8
$ws.Range("K106").Value = 6162.222
$ws.Range("L106").Value = 18149.8
$ws.Range("M106").Value = -5531.222
$ws.Range("N106").Value = -19411.8

$ws.Range("H107").Value = 697.6667
$ws.Range("I107").Value = 712.5333000000001
$ws.Range("J107").Value = 623.3333
$ws.Range("K107").Value = 712.5333000000001
$ws.Range("L107").Value = 623.3333
$ws.Range("M107").Value = 1207.4667
$ws.Range("N107").Value = -4463.3333

$ws.Range("H116").Value = 18856080
$ws.Range("I116").Value = 13216205
$ws.Range("J116").Value = 27785880
$ws.Range("K116").Value = 13216205
$ws.Range("L116").Value = 27785880
$ws.Range("M116").Value = -13212763
$ws.Range("N116").Value = -27792764

$ws.Range("H129").Value = 1566.0834
$ws.Range("I129").Value = 932.3333
$ws.Range("J129").Value = 1777.3334
$ws.Range("K129").Value = 2796.9999
$ws.Range("L129").Value = 5332.0002
$ws.Range("M129").Value = 2203.0001
$ws.Range("N129").Value = -15332.0002

$ws.Range("H137").Value = 3703.11
$ws.Range("I137").Value = 3035.2
$ws.Range("J137").Value = 3738.2632
$ws.Range("K137").Value = 9105.599999999999
$ws.Range("L137").Value = 11214.7896
$ws.Range("M137").Value = -6555.599999999999
$ws.Range("N137").Value = -16314.7896

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1975.5
$ws.Range("I2").Value = 1572
$ws.Range("J2").Value = 3993
$ws.Range("K2").Value = 1572
$ws.Range("L2").Value = 3993
$ws.Range("M2").Value = -1459
$ws.Range("N2").Value = -4219

$ws.Range("H116").Value = 1975.5
$ws.Range("I116").Value = 1572
$ws.Range("J116").Value = 3993
$ws.Range("K116").Value = 1572
$ws.Range("L116").Value = 3993
$ws.Range("M116").Value = 722
$ws.Range("N116").Value = -8581

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1975.5
$ws.Range("I3").Value = 1572
$ws.Range("J3").Value = 3993
$ws.Range("K3").Value = 1572
$ws.Range("L3").Value = 3993
$ws.Range("M3").Value = -1458
$ws.Range("N3").Value = -4221

$ws.Range("H105").Value = 1134.3948
$ws.Range("I105").Value = 1073.0769
$ws.Range("J105").Value = 1267.25
$ws.Range("K105").Value = 1073.0769
$ws.Range("L105").Value = 1267.25
$ws.Range("M105").Value = 673.9231
$ws.Range("N105").Value = -4761.25

$ws.Range("H107").Value = 50250910
$ws.Range("I107").Value = 301302.72
$ws.Range("J107").Value = 166800000
$ws.Range("K107").Value = 301302.72
$ws.Range("L107").Value = 166800000
$ws.Range("M107").Value = -299382.72
$ws.Range("N107").Value = -166803840

$ws.Range("H134").Value = 3285.5476
$ws.Range("I134").Value = 3317.5
$ws.Range("J134").Value = 3149.75
$ws.Range("K134").Value = 9952.5
$ws.Range("L134").Value = 9449.25
$ws.Range("M134").Value = -7417.5
$ws.Range("N134").Value = -14519.25

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H86").Value = 68770.766
$ws.Range("I86").Value = 107500.25
$ws.Range("J86").Value = 6803.6
$ws.Range("K86").Value = 107500.25
$ws.Range("L86").Value = 6803.6
$ws.Range("M86").Value = -106377.25
$ws.Range("N86").Value = -9049.6

$ws.Range("H89").Value = 68770.766
$ws.Range("I89").Value = 107500.25
$ws.Range("J89").Value = 6803.6
$ws.Range("K89").Value = 537501.25
$ws.Range("L89").Value = 34018
$ws.Range("M89").Value = -531885.25
$ws.Range("N89").Value = -45250

$ws.Range("I107").Value = 2820.3635
$ws.Range("J107").Value = 908.4
$ws.Range("K107").Value = 2820.3635
$ws.Range("L107").Value = 908.4
$ws.Range("M107").Value = -900.3634999999999
$ws.Range("N107").Value = -4748.4

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H68").Value = 2509.2
$ws.Range("I68").Value = 1912.5
$ws.Range("J68").Value = 2907
$ws.Range("K68").Value = 5737.5
$ws.Range("L68").Value = 8721
$ws.Range("M68").Value = -4926.5
$ws.Range("N68").Value = -10343

$ws.Range("H71").Value = 2509.2
$ws.Range("I71").Value = 1912.5
$ws.Range("J71").Value = 2907
$ws.Range("K71").Value = 17212.5
$ws.Range("L71").Value = 26163
$ws.Range("M71").Value = -13156.5
$ws.Range("N71").Value = -34275

$ws.Range("H132").Value = 1388.1538
$ws.Range("I132").Value = 1218.375
$ws.Range("J132").Value = 1659.8
$ws.Range("K132").Value = 10965.375
$ws.Range("L132").Value = 14938.2
$ws.Range("M132").Value = -8435.375
$ws.Range("N132").Value = -19998.2

$ws.Range("H140").Value = 1175.2106
$ws.Range("I140").Value = 1042.8823
$ws.Range("J140").Value = 2300
$ws.Range("K140").Value = 3128.6469
$ws.Range("L140").Value = 6900
$ws.Range("M140").Value = 2051.3531
$ws.Range("N140").Value = -17260

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 1670.9429
$ws.Range("I102").Value = 902.11536
$ws.Range("J102").Value = 3892
$ws.Range("K102").Value = 902.11536
$ws.Range("L102").Value = 3892
$ws.Range("M102").Value = 719.88464
$ws.Range("N102").Value = -7136

$ws.Range("H107").Value = 1332.1818
$ws.Range("I107").Value = 1119.875
$ws.Range("J107").Value = 1898.3334
$ws.Range("K107").Value = 1119.875
$ws.Range("L107").Value = 1898.3334
$ws.Range("M107").Value = 800.125
$ws.Range("N107").Value = -5738.3334

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 15627672
$ws.Range("I7").Value = 23811794
$ws.Range("J7").Value = 3436.7273
$ws.Range("K7").Value = 23811794
$ws.Range("L7").Value = 3436.7273
$ws.Range("M7").Value = -23811682
$ws.Range("N7").Value = -3660.7273

$ws.Range("H22").Value = 7500
$ws.Range("I22").Value = 7500
$ws.Range("J22").Value = 0
$ws.Range("K22").Value = 7500
$ws.Range("L22").Value = 0
$ws.Range("M22").Value = -7205
$ws.Range("N22").ClearContents()

$ws.Range("H27").Value = 7500
$ws.Range("I27").Value = 7500
$ws.Range("J27").Value = 0
$ws.Range("K27").Value = 7500
$ws.Range("L27").Value = 0
$ws.Range("M27").Value = -7393
$ws.Range("N27").ClearContents()

$ws.Range("H40").Value = 3137.9644
$ws.Range("I40").Value = 2882.9
$ws.Range("J40").Value = 3775.625
$ws.Range("K40").Value = 2882.9
$ws.Range("L40").Value = 3775.625
$ws.Range("M40").Value = -2746.9
$ws.Range("N40").Value = -4047.625

$ws.Range("H46").Value = 3536.2856
$ws.Range("I46").Value = 3459.8
$ws.Range("J46").Value = 3605.818
$ws.Range("K46").Value = 3459.8
$ws.Range("L46").Value = 3605.818
$ws.Range("M46").Value = -3271.8
$ws.Range("N46").Value = -3981.818

$ws.Range("H122").Value = 2848.087
$ws.Range("I122").Value = 2248.2727
$ws.Range("J122").Value = 4370.6924
$ws.Range("K122").Value = 6744.8181
$ws.Range("L122").Value = 13112.0772
$ws.Range("M122").Value = -4294.8181
$ws.Range("N122").Value = -18012.0772

$ws.Range("H126").Value = 15627672
$ws.Range("I126").Value = 23811794
$ws.Range("J126").Value = 3436.7273
$ws.Range("K126").Value = 71435382
$ws.Range("L126").Value = 10310.1819
$ws.Range("M126").Value = -71432912
$ws.Range("N126").Value = -15250.1819

$ws.Range("H136").Value = 3785.2856
$ws.Range("I136").Value = 1833.3334
$ws.Range("J136").Value = 5249.25
$ws.Range("K136").Value = 5500.0002
$ws.Range("L136").Value = 15747.75
$ws.Range("M136").Value = -2950.0002
$ws.Range("N136").Value = -20847.75

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 1901.75
$ws.Range("I122").Value = 1944.3684
$ws.Range("J122").Value = 1739.8
$ws.Range("K122").Value = 5833.1052
$ws.Range("L122").Value = 5219.4
$ws.Range("M122").Value = -3383.1052
$ws.Range("N122").Value = -10119.4

$ws.Range("H126").Value = 2082.1667
$ws.Range("I126").Value = 1898.875
$ws.Range("J126").Value = 2448.75
$ws.Range("K126").Value = 5696.625
$ws.Range("L126").Value = 7346.25
$ws.Range("M126").Value = -3226.625
$ws.Range("N126").Value = -12286.25

$ws.Range("H132").Value = 447915.12
$ws.Range("I132").Value = 529446.9399999999
$ws.Range("J132").Value = 5313.857
$ws.Range("K132").Value = 1588340.82
$ws.Range("L132").Value = 15941.571
$ws.Range("M132").Value = -1585810.82
$ws.Range("N132").Value = -21001.571
